$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C: "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Column C (previously audio file paths for the false/incorrect trial) now
# just records the current training phase name for both data rows.
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
